$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.486.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4665"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3739"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8893"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07968"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.857.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.437"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.604"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008968"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.513.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.190"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.057.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  +2.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.096"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.175"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08925"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.025"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7545"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.163"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.500"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.676"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.81%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01973"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.09%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05315"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.993"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.207"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5243"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1649"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.369"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4911"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.003"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.670"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06263"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.95%  "
